$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small edits to existing rows (balance tweaks) ---
$ws.Range("Q5").Value = 0.64
$ws.Range("Q9").Value = 0.13
$ws.Range("Q10").Value = 0.17
$ws.Range("S10").Value = -2

# --- Insert a new row at 11 for the new "Steyr AUG-Z Subcompact SWAT Keymod Top Rail" ---
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "steyr_aug_z_subcompact_swat_keymod_top_rail"
$ws.Range("B11").Value = "Steyr AUG-Z Subcompact SWAT Keymod Top Rail"
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 0.16
$ws.Range("R11").Value = -3
$ws.Range("S11").Value = -1
$ws.Range("Z11").Value = 1200
$ws.Range("AA11").Formula = "=P11-Q11*20-R11*0.8-S11*0.6-U11*5+V11*5+W11/300"

# --- Append new rows for Steyr AUG-Z sight rails after the existing data (row 35 is last, 36 blank sep) ---
$ws.Range("N36").Formula = "=C36-D36*20-E36*0.8-F36*0.6-H36*5+I36*5+J36/300"
$ws.Range("AA36").Formula = "=P36-Q36*20-R36*0.8-S36*0.6-U36*5+V36*5+W36/300"

$ws.Range("A37").Value = "steyr_aug_z_rail_front_sight"
$ws.Range("B37").Value = "Steyr AUG-Z Rail Front Sight"
$ws.Range("N37").Formula = "=C37-D37*20-E37*0.8-F37*0.6-H37*5+I37*5+J37/300"
$ws.Range("P37").Value = 3
$ws.Range("Q37").Value = 0.03
$ws.Range("Z37").Value = 0
$ws.Range("AA37").Formula = "=P37-Q37*20-R37*0.8-S37*0.6-U37*5+V37*5+W37/300"

$ws.Range("A38").Value = "steyr_aug_z_rail_front_sight_folded"
$ws.Range("B38").Value = "Steyr AUG-Z Rail Front Sight Folded"
$ws.Range("N38").Formula = "=C38-D38*20-E38*0.8-F38*0.6-H38*5+I38*5+J38/300"
$ws.Range("P38").Value = 2
$ws.Range("Q38").Value = 0.03
$ws.Range("Z38").Value = 0
$ws.Range("AA38").Formula = "=P38-Q38*20-R38*0.8-S38*0.6-U38*5+V38*5+W38/300"

$ws.Range("A39").Value = "steyr_aug_z_rail_rear_sight"
$ws.Range("B39").Value = "Steyr AUG-Z Rail Rear Sight "
$ws.Range("N39").Formula = "=C39-D39*20-E39*0.8-F39*0.6-H39*5+I39*5+J39/300"
$ws.Range("P39").Value = 3
$ws.Range("Q39").Value = 0.03
$ws.Range("Z39").Value = 0
$ws.Range("AA39").Formula = "=P39-Q39*20-R39*0.8-S39*0.6-U39*5+V39*5+W39/300"

$ws.Range("A40").Value = "steyr_aug_z_rail_rear_sight_folded"
$ws.Range("B40").Value = "Steyr AUG-Z Rail Rear Sight Folded"
$ws.Range("N40").Formula = "=C40-D40*20-E40*0.8-F40*0.6-H40*5+I40*5+J40/300"
$ws.Range("P40").Value = 2
$ws.Range("Q40").Value = 0.03
$ws.Range("Z40").Value = 0
$ws.Range("AA40").Formula = "=P40-Q40*20-R40*0.8-S40*0.6-U40*5+V40*5+W40/300"

# --- Column widths widened to fit new longer part names ---
$ws.Columns.Item(1).ColumnWidth = 30.86
$ws.Columns.Item(2).ColumnWidth = 41.57

# --- Selection moved ---
$ws.Range("Q6").Select()
